# Auto-generated edit script: updates odds values for rows 2-22
# per the betting-odds data refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("F2").Value = 1.62
$ws.Range("G2").Value = 1.64
$ws.Range("H2").Value = 5.8
$ws.Range("I2").Value = 6.2
$ws.Range("J2").Value = 4.5
$ws.Range("N2").Value = 5.5
$ws.Range("O2").Value = 1.21
$ws.Range("P2").Value = 2.5
$ws.Range("Q2").Value = 1.64
$ws.Range("R2").Value = 1.59
$ws.Range("S2").Value = 2.62
$ws.Range("T2").Value = 1.72
$ws.Range("U2").Value = 2.3
$ws.Range("V2").Value = 1.19
$ws.Range("W2").Value = 2.56
$ws.Range("Y2").Value = 26
$ws.Range("Z2").Value = 55
$ws.Range("AA2").Value = 150
$ws.Range("AC2").Value = 10.5
$ws.Range("AD2").Value = 22
$ws.Range("AE2").Value = 70
$ws.Range("AH2").Value = 18.5
$ws.Range("AI2").Value = 65
$ws.Range("AJ2").Value = 16
$ws.Range("AL2").Value = 29
$ws.Range("AM2").Value = 85
$ws.Range("AN2").Value = 6.8
$ws.Range("AO2").Value = 70
# Row 3
$ws.Range("F3").Value = 1.39
$ws.Range("G3").Value = 1.44
$ws.Range("P3").Value = 1.89
$ws.Range("S3").Value = 3.55
$ws.Range("W3").Value = 3.25
$ws.Range("AF3").Value = 9
# Row 4
$ws.Range("P4").Value = 2.22
$ws.Range("X4").Value = 19
$ws.Range("AJ4").Value = 85
# Row 5
$ws.Range("F5").Value = 1.71
$ws.Range("G5").Value = 1.73
$ws.Range("H5").Value = 5.2
$ws.Range("I5").Value = 5.4
$ws.Range("K5").Value = 4.5
$ws.Range("N5").Value = 5.6
$ws.Range("O5").Value = 1.2
$ws.Range("P5").Value = 2.56
$ws.Range("Q5").Value = 1.61
$ws.Range("R5").Value = 1.63
$ws.Range("S5").Value = 2.52
$ws.Range("T5").Value = 1.65
$ws.Range("U5").Value = 2.46
$ws.Range("V5").Value = 1.23
$ws.Range("W5").Value = 2.38
$ws.Range("X5").Value = 22
$ws.Range("Y5").Value = 25
$ws.Range("Z5").Value = 44
$ws.Range("AA5").Value = 120
$ws.Range("AC5").Value = 9.800000000000001
$ws.Range("AD5").Value = 19
$ws.Range("AF5").Value = 12
$ws.Range("AG5").Value = 9.800000000000001
$ws.Range("AH5").Value = 16.5
$ws.Range("AJ5").Value = 18
$ws.Range("AM5").Value = 70
$ws.Range("AN5").Value = 7.2
$ws.Range("AO5").Value = 48
# Row 6
$ws.Range("Q6").Value = 1.28
$ws.Range("X6").Value = 60
$ws.Range("AE6").Value = 85
$ws.Range("AF6").Value = 18
$ws.Range("AM6").Value = 70
# Row 7
$ws.Range("G7").Value = 2.14
$ws.Range("Q7").Value = 1.47
$ws.Range("S7").Value = 2.16
$ws.Range("W7").Value = 1.87
$ws.Range("X7").Value = 36
$ws.Range("Z7").Value = 980
$ws.Range("AA7").Value = 75
$ws.Range("AI7").Value = 42
# Row 8
$ws.Range("F8").Value = 2
$ws.Range("N8").Value = 4.6
$ws.Range("U8").Value = 2.52
$ws.Range("AA8").Value = 75
# Row 9
$ws.Range("G9").Value = 1.5
$ws.Range("S9").Value = 1.98
$ws.Range("T9").Value = 1.64
$ws.Range("AA9").Value = 220
$ws.Range("AC9").Value = 17
$ws.Range("AO9").Value = 80
# Row 10
$ws.Range("G10").Value = 2.26
$ws.Range("M10").Value = 1.1
$ws.Range("U10").Value = 1.75
$ws.Range("W10").Value = 1.79
# Row 11
$ws.Range("F11").Value = 2.74
$ws.Range("G11").Value = 2.98
$ws.Range("H11").Value = 2.58
$ws.Range("I11").Value = 2.78
$ws.Range("J11").Value = 3.5
$ws.Range("K11").Value = 3.75
$ws.Range("N11").Value = 3.9
$ws.Range("P11").Value = 2.02
$ws.Range("Q11").Value = 1.86
$ws.Range("R11").Value = 1.4
$ws.Range("S11").Value = 3.15
$ws.Range("T11").Value = 1.69
$ws.Range("U11").Value = 2.24
$ws.Range("X11").Value = 19
$ws.Range("Y11").Value = 13.5
$ws.Range("AB11").Value = 14.5
$ws.Range("AC11").Value = 9.199999999999999
$ws.Range("AD11").Value = 14.5
# Row 12
$ws.Range("L12").Value = 1.48
$ws.Range("M12").Value = 1.08
$ws.Range("V12").Value = 1.18
# Row 13
$ws.Range("F13").Value = 1.91
$ws.Range("G13").Value = 1.92
$ws.Range("H13").Value = 4.4
$ws.Range("I13").Value = 4.6
$ws.Range("L13").Value = 1.34
$ws.Range("M13").Value = 1.05
$ws.Range("P13").Value = 2.34
$ws.Range("Q13").Value = 1.72
$ws.Range("U13").Value = 2.38
$ws.Range("V13").Value = 1.28
$ws.Range("W13").Value = 2.08
$ws.Range("X13").Value = 19
$ws.Range("Z13").Value = 34
$ws.Range("AA13").Value = 90
$ws.Range("AB13").Value = 11.5
$ws.Range("AD13").Value = 16.5
$ws.Range("AF13").Value = 12.5
$ws.Range("AH13").Value = 16
$ws.Range("AI13").Value = 48
$ws.Range("AN13").Value = 9.800000000000001
# Row 14
$ws.Range("F14").Value = 3.35
$ws.Range("H14").Value = 2.18
$ws.Range("L14").Value = 1.3
$ws.Range("P14").Value = 2.5
$ws.Range("Q14").Value = 1.63
$ws.Range("S14").Value = 2.56
$ws.Range("T14").Value = 1.57
$ws.Range("V14").Value = 1.81
$ws.Range("W14").Value = 1.4
$ws.Range("AO14").Value = 11.5
# Row 15
$ws.Range("P15").Value = 3.25
$ws.Range("R15").Value = 1.93
$ws.Range("W15").Value = 1.67
$ws.Range("AC15").Value = 10.5
$ws.Range("AH15").Value = 13
# Row 16
$ws.Range("H16").Value = 15
$ws.Range("I16").Value = 16
$ws.Range("P16").Value = 3.45
$ws.Range("Q16").Value = 1.38
$ws.Range("R16").Value = 2
$ws.Range("S16").Value = 1.94
$ws.Range("T16").Value = 1.89
$ws.Range("U16").Value = 2.04
$ws.Range("V16").Value = 1.06
$ws.Range("Y16").Value = 70
$ws.Range("AA16").Value = 590
$ws.Range("AC16").Value = 18
$ws.Range("AG16").Value = 12
$ws.Range("AH16").Value = 30
$ws.Range("AO16").Value = 160
# Row 17
$ws.Range("J17").Value = 6.4
$ws.Range("K17").Value = 6.6
$ws.Range("P17").Value = 2.58
$ws.Range("R17").Value = 1.62
$ws.Range("S17").Value = 2.5
$ws.Range("U17").Value = 1.89
$ws.Range("AB17").Value = 9.800000000000001
$ws.Range("AJ17").Value = 9.800000000000001
# Row 18
$ws.Range("H18").Value = 2.6
$ws.Range("I18").Value = 2.62
$ws.Range("L18").Value = 1.37
$ws.Range("Q18").Value = 1.85
$ws.Range("X18").Value = 16
# Row 19
$ws.Range("H19").Value = 2.46
$ws.Range("I19").Value = 2.48
$ws.Range("O19").Value = 1.37
$ws.Range("T19").Value = 1.83
$ws.Range("V19").Value = 1.67
$ws.Range("AA19").Value = 34
$ws.Range("AI19").Value = 42
$ws.Range("AM19").Value = 100
$ws.Range("AO19").Value = 23
# Row 20
$ws.Range("G20").Value = 4.9
$ws.Range("J20").Value = 3.5
$ws.Range("N20").Value = 3.2
$ws.Range("P20").Value = 1.78
$ws.Range("R20").Value = 1.29
$ws.Range("S20").Value = 4
$ws.Range("U20").Value = 1.92
$ws.Range("V20").Value = 2.02
$ws.Range("W20").Value = 1.25
$ws.Range("Y20").Value = 9.199999999999999
# Row 21
$ws.Range("F21").Value = 3.05
$ws.Range("G21").Value = 3.15
$ws.Range("K21").Value = 3.45
$ws.Range("N21").Value = 3.3
$ws.Range("W21").Value = 1.46
# Row 22
$ws.Range("M22").Value = 1.08
$ws.Range("O22").Value = 1.42
$ws.Range("Q22").Value = 2.08
$ws.Range("S22").Value = 3.75
